$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New profile entry (row 12): name + github handle, same pattern as the
# existing rows. New shared strings get appended automatically.
$ws.Range("A12").Value = "Shivankshi Khandelwal"
$ws.Range("B12").Value = "shivankshi13"

# Columns A (names) and B (handles) get wider, fixed widths so the new,
# longer entries aren't truncated; the rest of the sheet keeps the
# worksheet's default column width.
$ws.Columns.Item(1).ColumnWidth = 23.7
$ws.Columns.Item(2).ColumnWidth = 18.85

# Leave the selection where the next entry would be typed in.
$ws.Range("B13").Select()
